# Insert a new data row at row 291 (pushes existing rows 291..365 down to 292..366)
# and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("291:291").Insert()

$ws.Range("A291").Value = 6
$ws.Range("B291").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C291").Value = "Metropolitana"
$ws.Range("D291").Value = 45204
$ws.Range("E291").Value = 13
$ws.Range("F291").Value = 100112001
$ws.Range("G291").Value = "Berenjena"
$ws.Range("H291").Value = "Sin especificar"
$ws.Range("I291").Value = "Primera"
$ws.Range("J291").Value = 250
$ws.Range("K291").Value = 6000
$ws.Range("L291").Value = 7000
$ws.Range("M291").Value = 6480
$ws.Range("N291").Value = "$/caja 50 unidades"
$ws.Range("O291").Value = "Región de Arica y Parinacota"
$ws.Range("P291").Value = 130
$ws.Range("Q291").Value = 50
$ws.Range("R291").Value = "Hortaliza"
